# Auto-generated script to update Famfrit Profits market data values
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H10").Value = 500
$ws.Range("H38").Value = 4714.769
$ws.Range("J38").Value = 5033.6665
$ws.Range("L38").Value = 15100.9995
$ws.Range("N38").Value = -15844.9995
$ws.Range("K62").Value = 2125
$ws.Range("N62").Value = -7987
$ws.Range("L62").Value = 6739
$ws.Range("M62").Value = -1501
$ws.Range("J62").Value = 6739
$ws.Range("H62").Value = 5816.2
$ws.Range("I62").Value = 2125
$ws.Range("H65").Value = 5816.2
$ws.Range("I65").Value = 2125
$ws.Range("J65").Value = 6739
$ws.Range("M65").Value = -7505
$ws.Range("K65").Value = 10625
$ws.Range("N65").Value = -39935
$ws.Range("L65").Value = 33695
$ws.Range("H74").Value = 4529.1665
$ws.Range("I74").Value = 2599.25
$ws.Range("M74").Value = -1663.25
$ws.Range("K74").Value = 2599.25
$ws.Range("H77").Value = 4529.1665
$ws.Range("K77").Value = 12996.25
$ws.Range("M77").Value = -8316.25
$ws.Range("I77").Value = 2599.25
$ws.Range("I80").Value = 4250
$ws.Range("H80").Value = 4768.5386
$ws.Range("K80").Value = 12750
$ws.Range("M80").Value = -11752
$ws.Range("I83").Value = 4250
$ws.Range("M83").Value = -33258
$ws.Range("H83").Value = 4768.5386
$ws.Range("K83").Value = 38250
$ws.Range("K86").Value = 2314.5625
$ws.Range("I86").Value = 2314.5625
$ws.Range("M86").Value = -1191.5625
$ws.Range("H86").Value = 3163.1292
$ws.Range("M88").Value = -6844
$ws.Range("J88").Value = 6962.3335
$ws.Range("K88").Value = 7250
$ws.Range("I88").Value = 7250
$ws.Range("H88").Value = 7077.4
$ws.Range("N88").Value = -7774.3335
$ws.Range("L88").Value = 6962.3335
$ws.Range("M89").Value = -5956.8125
$ws.Range("I89").Value = 2314.5625
$ws.Range("K89").Value = 11572.8125
$ws.Range("H89").Value = 3163.1292
$ws.Range("N91").Value = -9770.333500000001
$ws.Range("J91").Value = 6962.3335
$ws.Range("I91").Value = 7250
$ws.Range("M91").Value = -5846
$ws.Range("H91").Value = 7077.4
$ws.Range("K91").Value = 7250
$ws.Range("L91").Value = 6962.3335
$ws.Range("M96").Value = -7981
$ws.Range("H96").Value = 4425
$ws.Range("I96").Value = 3118
$ws.Range("K96").Value = 9354
$ws.Range("K98").Value = 1562
$ws.Range("M98").Value = -64
$ws.Range("H98").Value = 1562
$ws.Range("I98").Value = 1562
$ws.Range("H100").Value = 2990.2
$ws.Range("M100").Value = -959.6666
$ws.Range("K100").Value = 1500.6666
$ws.Range("I100").Value = 1500.6666
$ws.Range("K111").Value = 50007678
$ws.Range("H111").Value = 13335706
$ws.Range("I111").Value = 16669226
$ws.Range("J111").Value = 1626
$ws.Range("L111").Value = 4878
$ws.Range("M111").Value = -50004611
$ws.Range("N111").Value = -11012
$ws.Range("M113").Value = 556.2222000000002
$ws.Range("H113").Value = 4706
$ws.Range("K113").Value = 2697.7778
$ws.Range("I113").Value = 2697.7778
$ws.Range("L116").Value = 3999.6667
$ws.Range("I116").Value = 2756.5557
$ws.Range("J116").Value = 3999.6667
$ws.Range("M116").Value = 685.4443000000001
$ws.Range("H116").Value = 3067.3333
$ws.Range("K116").Value = 2756.5557
$ws.Range("N116").Value = -10883.6667
$ws.Range("K122").Value = 4686
$ws.Range("M122").Value = -2236
$ws.Range("I122").Value = 1562
$ws.Range("H122").Value = 1562
$ws.Range("M138").Value = 1125.6667
$ws.Range("L138").Value = 55582302
$ws.Range("K138").Value = 4014.3333
$ws.Range("H138").Value = 12352069
$ws.Range("N138").Value = -55592582
$ws.Range("J138").Value = 18527434
$ws.Range("I138").Value = 1338.1111
$ws.Range("H141").Value = 1842.0834
$ws.Range("J141").Value = 3103.3333
$ws.Range("K141").Value = 4985.7144
$ws.Range("L141").Value = 9309.999899999999
$ws.Range("M141").Value = 194.2856000000002
$ws.Range("N141").Value = -19669.9999
$ws.Range("I141").Value = 1661.9048

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -1059
$ws.Range("H2").Value = 890.38464
$ws.Range("I2").Value = 907.6
$ws.Range("J2").Value = 833
$ws.Range("M2").Value = -794.6
$ws.Range("K2").Value = 907.6
$ws.Range("L2").Value = 833
$ws.Range("I10").Value = 6701.6665
$ws.Range("K10").Value = 6701.6665
$ws.Range("M10").Value = -6531.6665
$ws.Range("H10").Value = 6701.6665
$ws.Range("I32").Value = 9724.294
$ws.Range("K32").Value = 9724.294
$ws.Range("H32").Value = 10072.944
$ws.Range("M32").Value = -9437.294
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10828
$ws.Range("H41").Value = 6749.5
$ws.Range("H61").Value = 45460264
$ws.Range("M61").Value = -55560112
$ws.Range("K61").Value = 55560324
$ws.Range("I61").Value = 55560324
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10424
$ws.Range("N110").Value = -8504.833500000001
$ws.Range("K110").Value = 16100.218
$ws.Range("H110").Value = 13682.552
$ws.Range("I110").Value = 16100.218
$ws.Range("L110").Value = 4414.8335
$ws.Range("M110").Value = -14055.218
$ws.Range("J110").Value = 4414.8335
$ws.Range("J112").Value = 24999
$ws.Range("N112").Value = -27953
$ws.Range("L112").Value = 24999
$ws.Range("H112").Value = 24999
$ws.Range("L116").Value = 833
$ws.Range("I116").Value = 907.6
$ws.Range("J116").Value = 833
$ws.Range("M116").Value = 1386.4
$ws.Range("H116").Value = 890.38464
$ws.Range("K116").Value = 907.6
$ws.Range("N116").Value = -5421
$ws.Range("K122").Value = 7381.875
$ws.Range("M122").Value = -4931.875
$ws.Range("I122").Value = 2460.625
$ws.Range("H122").Value = 3249.4075
$ws.Range("K132").Value = 34434.783
$ws.Range("H132").Value = 34540724
$ws.Range("I132").Value = 11478.261
$ws.Range("M132").Value = -31904.783
$ws.Range("K136").Value = 166680972
$ws.Range("I136").Value = 55560324
$ws.Range("M136").Value = -166678422
$ws.Range("N136").Value = -35100
$ws.Range("J136").Value = 10000
$ws.Range("H136").Value = 45460264
$ws.Range("L136").Value = 30000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 833
$ws.Range("H3").Value = 890.38464
$ws.Range("K3").Value = 907.6
$ws.Range("I3").Value = 907.6
$ws.Range("M3").Value = -793.6
$ws.Range("L3").Value = 833
$ws.Range("N3").Value = -1061
$ws.Range("K20").Value = 2511.8
$ws.Range("L20").Value = 2027.3334
$ws.Range("M20").Value = -2264.8
$ws.Range("I20").Value = 2511.8
$ws.Range("J20").Value = 2027.3334
$ws.Range("H20").Value = 2247.5454
$ws.Range("N20").Value = -2521.3334
$ws.Range("K86").Value = 21585.5
$ws.Range("I86").Value = 21585.5
$ws.Range("M86").Value = -20462.5
$ws.Range("H86").Value = 81057
$ws.Range("M89").Value = -102311.5
$ws.Range("I89").Value = 21585.5
$ws.Range("K89").Value = 107927.5
$ws.Range("H89").Value = 81057
$ws.Range("H107").Value = 3425.8
$ws.Range("N107").Value = -10586.5
$ws.Range("L107").Value = 6746.5
$ws.Range("I107").Value = 2914.923
$ws.Range("K107").Value = 2914.923
$ws.Range("M107").Value = -994.9229999999998
$ws.Range("J107").Value = 6746.5
$ws.Range("I134").Value = 2627.111
$ws.Range("K134").Value = 7881.333
$ws.Range("M134").Value = -5346.333
$ws.Range("H134").Value = 3131.88

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N31").Value = -7113.231
$ws.Range("K31").Value = 2283.7026
$ws.Range("H31").Value = 3385.98
$ws.Range("L31").Value = 6523.231
$ws.Range("M31").Value = -1988.7026
$ws.Range("J31").Value = 6523.231
$ws.Range("I31").Value = 2283.7026
$ws.Range("N34").Value = -6927.231
$ws.Range("M34").Value = -2081.7026
$ws.Range("I34").Value = 2283.7026
$ws.Range("J34").Value = 6523.231
$ws.Range("K34").Value = 2283.7026
$ws.Range("L34").Value = 6523.231
$ws.Range("H34").Value = 3385.98
$ws.Range("N57").ClearContents()
$ws.Range("J57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("L67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -12499.556
$ws.Range("N99").ClearContents()
$ws.Range("H99").Value = 13997.556
$ws.Range("K99").Value = 13997.556
$ws.Range("I99").Value = 13997.556
$ws.Range("J99").Value = 0
$ws.Range("H107").Value = 1488.8334
$ws.Range("N107").Value = -7421
$ws.Range("L107").Value = 3581
$ws.Range("I107").Value = 627.35297
$ws.Range("K107").Value = 627.35297
$ws.Range("M107").Value = 1292.64703
$ws.Range("J107").Value = 3581
$ws.Range("K122").Value = 4153.8333
$ws.Range("M122").Value = -1703.8333
$ws.Range("I122").Value = 1384.6111
$ws.Range("J122").Value = 4900
$ws.Range("N122").Value = -19600
$ws.Range("L122").Value = 14700
$ws.Range("H122").Value = 1569.6316
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 41992.66800000001
$ws.Range("H126").Value = 13997.556
$ws.Range("L126").Value = 0
$ws.Range("I126").Value = 13997.556
$ws.Range("M126").Value = -39522.66800000001
$ws.Range("N126").ClearContents()
$ws.Range("H141").Value = 264959.38
$ws.Range("J141").Value = 264959.38
$ws.Range("L141").Value = 264959.38
$ws.Range("N141").Value = -275319.38

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1219.9584
$ws.Range("N5").Value = -8046.7145
$ws.Range("L5").Value = 7822.7145
$ws.Range("J5").Value = 2607.5715
$ws.Range("N12").Value = -2530.3
$ws.Range("J12").Value = 728.1
$ws.Range("K12").Value = 1074
$ws.Range("M12").Value = -901
$ws.Range("H12").Value = 622.3570999999999
$ws.Range("I12").Value = 358
$ws.Range("L12").Value = 2184.3
$ws.Range("H37").Value = 198936.6
$ws.Range("L37").Value = 596809.8
$ws.Range("N37").Value = -597033.8
$ws.Range("J37").Value = 198936.6
$ws.Range("L135").Value = 23468.1435
$ws.Range("N135").Value = -28538.1435
$ws.Range("J135").Value = 2607.5715
$ws.Range("H135").Value = 1219.9584

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L49").Value = 30747.5
$ws.Range("N49").Value = -31115.5
$ws.Range("H49").Value = 30099
$ws.Range("J49").Value = 30747.5
$ws.Range("H107").Value = 676.9583
$ws.Range("N107").Value = -4721.625
$ws.Range("L107").Value = 881.625
$ws.Range("I107").Value = 574.625
$ws.Range("K107").Value = 574.625
$ws.Range("M107").Value = 1345.375
$ws.Range("J107").Value = 881.625
$ws.Range("J132").Value = 2542.5715
$ws.Range("K132").Value = 5798.7693
$ws.Range("H132").Value = 2146.3
$ws.Range("I132").Value = 1932.9231
$ws.Range("M132").Value = -3268.7693
$ws.Range("N132").Value = -12687.7145
$ws.Range("L132").Value = 7627.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L16").Value = 1100
$ws.Range("N16").Value = -1440
$ws.Range("M16").Value = -830
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1100
$ws.Range("H16").Value = 1066.6666
$ws.Range("K16").Value = 1000
$ws.Range("H63").Value = 66749.5
$ws.Range("N63").Value = -65497
$ws.Range("L63").Value = 63999
$ws.Range("J63").Value = 63999
$ws.Range("L66").Value = 191997
$ws.Range("J66").Value = 63999
$ws.Range("H66").Value = 66749.5
$ws.Range("N66").Value = -199485
$ws.Range("N81").Value = -60996
$ws.Range("H81").Value = 74000
$ws.Range("J81").Value = 59000
$ws.Range("L81").Value = 59000
$ws.Range("H84").Value = 74000
$ws.Range("L84").Value = 177000
$ws.Range("N84").Value = -186984
$ws.Range("J84").Value = 59000
$ws.Range("K122").Value = 8294.5386
$ws.Range("M122").Value = -5844.5386
$ws.Range("I122").Value = 2764.8462
$ws.Range("H122").Value = 4460.8276
$ws.Range("K136").Value = 5673.272999999999
$ws.Range("I136").Value = 1891.091
$ws.Range("M136").Value = -3123.272999999999
$ws.Range("H136").Value = 2487.5625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 750000
$ws.Range("I2").Value = 750000
$ws.Range("M2").Value = -749888
$ws.Range("K2").Value = 750000
$ws.Range("H15").Value = 10250
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("L70").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("H92").Value = 10000
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("N103").Value = -19496.4
$ws.Range("H103").Value = 17152.4
$ws.Range("J103").Value = 17152.4
$ws.Range("L103").Value = 17152.4
$ws.Range("H107").Value = 429.4
$ws.Range("I107").Value = 463.42856
$ws.Range("K107").Value = 1390.28568
$ws.Range("M107").Value = 529.71432
$ws.Range("M113").Value = 607.6000000000001
$ws.Range("H113").Value = 908.45
$ws.Range("K113").Value = 1562.4
$ws.Range("I113").Value = 520.8
$ws.Range("K122").Value = 128896.125
$ws.Range("M122").Value = -126446.125
$ws.Range("I122").Value = 42965.375
$ws.Range("H122").Value = 34870.434
$ws.Range("K136").Value = 3196.5
$ws.Range("I136").Value = 1065.5
$ws.Range("M136").Value = -646.5
$ws.Range("H136").Value = 1871.4584
